$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fit" column (G) for the Jeans rows (26-37) was "Baggy"; rename it to
# "Pantalón Baggy" so the garment type reads "Pantalón Baggy" instead of
# plain "Baggy". Writing the same text to every cell in the range lets the
# engine drop the now-unused "Baggy" shared string and append the new one,
# matching the sharedStrings.xml renumbering in the diff.
$ws.Range("G26:G37").Value2 = "Pantalón Baggy"

# Widen column G (Fit) and column M (Código Barras) so the longer label and
# barcode values are fully visible.
$ws.Columns.Item(7).ColumnWidth = 19.91666666666667
$ws.Columns.Item(13).ColumnWidth = 47.41666666666667

# Move the active selection to G27, reflecting where the editor was working.
$ws.Range("G27").Select()
